$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 02:23:48"

$ws1.Range("A6").Value = "02:23:48"
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 35

$ws1.Range("A7").Value = "02:23:48"
$ws1.Range("B7").Value = "03:55"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 92

$ws1.Range("A8").Value = "02:23:48"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 98

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 02:23:48"

$ws2.Range("A6").Value = "02:23:48"
$ws2.Range("B6").Value = "02:58"
$ws2.Range("D6").Value = 35

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 02:23:48"
